# Ran analysis with new .csv files -- update diff/CI/pval results (rows 2-17, cols C:G)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,5

$data[0,0] = 0.04
$data[0,1] = -1.10493476989594
$data[0,2] = 1.18493476989594
$data[0,3] = 0.999954040003884
$data[0,4] = 0.999954040003884

$data[1,0] = 0.04
$data[1,1] = -1.10493476989594
$data[1,2] = 1.18493476989594
$data[1,3] = 0.999954024347963
$data[1,4] = 0.999954040003884

$data[2,0] = -0.05
$data[2,1] = -1.19493476989594
$data[2,2] = 1.09493476989594
$data[2,3] = 0.999888453518731
$data[2,4] = 0.999954040003884

$data[3,0] = 0.56
$data[3,1] = -0.584934769895937
$data[3,2] = 1.70493476989594
$data[3,3] = 0.546841751127347
$data[3,4] = 0.972163113115283

$data[4,0] = 7.29
$data[4,1] = -3.86765921809915
$data[4,2] = 18.4476592180991
$data[4,3] = 0.296228491253824
$data[4,4] = 0.789942643343531

$data[5,0] = 1.75
$data[5,1] = -9.40765921809915
$data[5,2] = 12.9076592180991
$data[5,3] = 0.98423965178173
$data[5,4] = 0.999954040003884

$data[6,0] = 2.91
$data[6,1] = -8.24765921809915
$data[6,2] = 14.0676592180991
$data[6,3] = 0.909667303939179
$data[6,4] = 0.999954040003884

$data[7,0] = 6.06
$data[7,1] = -5.09765921809915
$data[7,2] = 17.2176592180991
$data[7,3] = 0.455723732455904
$data[7,4] = 0.911447464911807

$data[8,0] = 3.6
$data[8,1] = -13.8966719397445
$data[8,2] = 21.0966719397445
$data[8,3] = 0.958643776734753
$data[8,4] = 0.999954040003884

$data[9,0] = 10.26
$data[9,1] = -7.23667193974447
$data[9,2] = 27.7566719397445
$data[9,3] = 0.388376644157544
$data[9,4] = 0.887718043788671

$data[10,0] = 13.39
$data[10,1] = -4.10667193974447
$data[10,2] = 30.8866719397445
$data[10,3] = 0.178078800092171
$data[10,4] = 0.569852160294948

$data[11,0] = 4.21
$data[11,1] = -13.2866719397445
$data[11,2] = 21.7066719397445
$data[11,3] = 0.93016383182917
$data[11,4] = 0.999954040003884

$data[12,0] = 12.24
$data[12,1] = -3.46541705981242
$data[12,2] = 27.9454170598124
$data[12,3] = 0.166259632852575
$data[12,4] = 0.569852160294948

$data[13,0] = 27.74
$data[13,1] = 12.0345829401876
$data[13,2] = 43.4454170598124
$data[13,3] = 0.000206099036646368
$data[13,4] = 0.00109919486211396

$data[14,0] = 52.24
$data[14,1] = 36.5345829401876
$data[14,2] = 67.9454170598124
$data[14,3] = 0.0000000000274595901572638
$data[14,4] = 0.000000000439353442516222

$data[15,0] = 43.18
$data[15,1] = 27.4745829401876
$data[15,2] = 58.8854170598124
$data[15,3] = 0.0000000245471284410215
$data[15,4] = 0.000000196377027528172

$ws.Range("C2:G17").Value = $data
